$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($rng, $val)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.860.12'
Set-TextValue $ws.Range("E2") '  -1.79%  '
Set-TextValue $ws.Range("D3") '1.801.71'
Set-TextValue $ws.Range("E3") '  -1.30%  '
Set-TextValue $ws.Range("E4") '  -0.08%  '
Set-TextValue $ws.Range("D5") '309.15'
Set-TextValue $ws.Range("E5") '  -1.73%  '
Set-TextValue $ws.Range("E6") '  -0.05%  '
Set-TextValue $ws.Range("D7") '0.4667'
Set-TextValue $ws.Range("E7") '  +4.34%  '
Set-TextValue $ws.Range("D8") '0.3698'
Set-TextValue $ws.Range("E8") '  -2.09%  '
Set-TextValue $ws.Range("D9") '0.07380'
Set-TextValue $ws.Range("E9") '  -0.89%  '
Set-TextValue $ws.Range("D10") '0.8698'
Set-TextValue $ws.Range("E10") '  -2.28%  '
Set-TextValue $ws.Range("E11") '  -3.11%  '
Set-TextValue $ws.Range("D12") '1.748.25'
Set-TextValue $ws.Range("E12") '  -4.26%  '
Set-TextValue $ws.Range("D13") '5.357'
Set-TextValue $ws.Range("E13") '  -2.03%  '
Set-TextValue $ws.Range("D14") '92.31'
Set-TextValue $ws.Range("E14") '  -1.42%  '
Set-TextValue $ws.Range("D15") '6.495'
Set-TextValue $ws.Range("E15") '  -3.74%  '
Set-TextValue $ws.Range("E16") '  -1.48%  '
Set-TextValue $ws.Range("D17") '1.001'
Set-TextValue $ws.Range("E17") '  -0.03%  '
Set-TextValue $ws.Range("D18") '0.000008705'
Set-TextValue $ws.Range("E18") '  -1.00%  '
Set-TextValue $ws.Range("D19") '1.001'
Set-TextValue $ws.Range("E19") '  -0.02%  '
Set-TextValue $ws.Range("D20") '14.69'
Set-TextValue $ws.Range("E20") '  -3.07%  '
Set-TextValue $ws.Range("D21") '26.848.36'
Set-TextValue $ws.Range("E21") '  -1.84%  '
Set-TextValue $ws.Range("D22") '5.297'
Set-TextValue $ws.Range("E22") '  -1.85%  '
Set-TextValue $ws.Range("D23") '10.61'
Set-TextValue $ws.Range("E23") '  -3.38%  '
Set-TextValue $ws.Range("D24") '2.062.29'
Set-TextValue $ws.Range("E24") '  +0.42%  '
Set-TextValue $ws.Range("D25") '1.902'
Set-TextValue $ws.Range("E25") '  -3.46%  '
Set-TextValue $ws.Range("D26") '151.57'
Set-TextValue $ws.Range("E26") '  +0.06%  '
Set-TextValue $ws.Range("E27") '  -1.90%  '
Set-TextValue $ws.Range("D28") '2.138'
Set-TextValue $ws.Range("E28") '  -8.42%  '
Set-TextValue $ws.Range("D29") '5.261'
Set-TextValue $ws.Range("E29") '  -2.43%  '
Set-TextValue $ws.Range("D30") '115.80'
Set-TextValue $ws.Range("E30") '  -1.78%  '
Set-TextValue $ws.Range("D31") '0.08931'
Set-TextValue $ws.Range("E31") '  +0.53%  '
Set-TextValue $ws.Range("D32") '0.7563'
Set-TextValue $ws.Range("E32") '  -4.95%  '
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D33") '1.150'
Set-TextValue $ws.Range("E33") '  -4.40%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D34") '2.920'
Set-TextValue $ws.Range("E34") '  +0.03%  '
Set-TextValue $ws.Range("D35") '4.450'
Set-TextValue $ws.Range("E35") '  -3.54%  '
Set-TextValue $ws.Range("D36") '1.000'
Set-TextValue $ws.Range("E36") '  -0.04%  '
Set-TextValue $ws.Range("E37") '  -1.22%  '
Set-TextValue $ws.Range("D38") '0.01955'
Set-TextValue $ws.Range("E38") '  -1.64%  '
Set-TextValue $ws.Range("D39") '0.05246'
Set-TextValue $ws.Range("E39") '  -1.31%  '
Set-TextValue $ws.Range("E40") '  +1.84%  '
Set-TextValue $ws.Range("D41") '7.243'
Set-TextValue $ws.Range("E41") '  -0.80%  '
Set-TextValue $ws.Range("D42") '2.384'
Set-TextValue $ws.Range("E42") '  +2.13%  '
Set-TextValue $ws.Range("D43") '0.5279'
Set-TextValue $ws.Range("E43") '  -1.62%  '
Set-TextValue $ws.Range("E44") '  -3.59%  '
Set-TextValue $ws.Range("D45") '8.487'
Set-TextValue $ws.Range("E45") '  -2.21%  '
Set-TextValue $ws.Range("D46") '0.4994'
Set-TextValue $ws.Range("E46") '  -2.03%  '
Set-TextValue $ws.Range("D47") '10.28'
Set-TextValue $ws.Range("E47") '  -3.20%  '
Set-TextValue $ws.Range("D48") '104.13'
Set-TextValue $ws.Range("E48") '  -1.10%  '
Set-TextValue $ws.Range("D49") '1.0000'
Set-TextValue $ws.Range("E49") '  -0.03%  '
Set-TextValue $ws.Range("D50") '1.663'
Set-TextValue $ws.Range("E50") '  -2.05%  '
Set-TextValue $ws.Range("D51") '0.06288'
Set-TextValue $ws.Range("E51") '  -1.91%  '
